$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Step 1: Remove the standalone "Meta description" paragraph that
# originally followed the page title (it gets dropped from the top of
# the document entirely).
# ------------------------------------------------------------------
$removed = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Meta description*") {
        $p.Range.Delete()
        $removed = $true
        break
    }
}
Write-Host "Removed meta description paragraph: " $removed

# ------------------------------------------------------------------
# Step 2: Find the paragraph that holds the old image-generation
# prompt text (the last paragraph in the document) and replace it with
# two new paragraphs:
#   1) a bold heading-style line: "Play Firebird Double 27 for Free | Review"
#   2) an italic line with the meta-description copy that used to live
#      near the top of the document: "Read our review of Firebird
#      Double 27 and play it for free. This slot game offers a chance
#      to win significant jackpots to all players."
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Please create a feature image*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $xmlNew = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
              "<w:r/><w:r><w:rPr><w:b/></w:rPr>" +
              "<w:t>Play Firebird Double 27 for Free | Review</w:t></w:r></w:p>" +
              "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
              "<w:r/><w:r><w:rPr><w:i/></w:rPr>" +
              "<w:t>Read our review of Firebird Double 27 and play it for free. This slot game offers a chance to win significant jackpots to all players.</w:t></w:r></w:p>"
    $rng.InsertXML($xmlNew)
    Write-Host "Replaced final paragraph with title + description paragraphs"
} else {
    Write-Host "WARNING: target paragraph (image prompt) not found!"
}
